$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("A2").Value = "DEL-0179"
$ws.Range("B2").Value = "DSR-0248"
$ws.Range("C2").Value = "Arshi Electronics"
$ws.Range("D2").Value = "Loxmicole"
$ws.Range("E2").Value = "Md Abdullah Al Mamun"
$ws.Range("G2").Value = "GO"
$ws.Range("I2").Value = "Md Abdullah Al Mamun"
$ws.Range("J2").Value = 1726359224
$ws.Range("K2").Value = "Natore"
$ws.Range("L2").Value = "Baraigram"
$ws.Range("M2").Value = "ZSO-0022"
$ws.Range("N2").Value = "Laxmicole, Baraigram, Natore."
$ws.Range("P2").Value = 1726359224
$ws.Range("Q2").Value = "C"
$ws.Range("R2").Value = "Rural"
$ws.Range("S2").Value = "bKash"
$ws.Range("T2").Value = 1726359224

$ws.Range("A3").Value = "DEL-0179"
$ws.Range("B3").Value = "DSR-0349"
$ws.Range("C3").Value = "Ma Telecom "
$ws.Range("D3").Value = "Bagha"
$ws.Range("E3").Value = "Md Babu Hosen"
$ws.Range("G3").Value = "GO"
$ws.Range("I3").Value = "Md Babu Hosen"
$ws.Range("J3").Value = 1740418484
$ws.Range("K3").Value = "Rajshahi"
$ws.Range("L3").Value = "Bagha"
$ws.Range("M3").Value = "ZSO-0022"
$ws.Range("N3").Value = "Arani Bazar, Bagha, Rajshahi."
$ws.Range("P3").Value = 1740418484
$ws.Range("Q3").Value = "C"
$ws.Range("R3").Value = "Rural"
$ws.Range("S3").Value = "bKash"
$ws.Range("T3").Value = 1740418484

$ws.Range("A4").Value = "DEL-0179"
$ws.Range("B4").Value = "DSR-0350"
$ws.Range("C4").Value = "Nabinogor Bohumukhi Somobai Somiti"
$ws.Range("D4").Value = "Lalpur"
$ws.Range("E4").Value = "Md Biplob Hossain"
$ws.Range("G4").Value = "GO"
$ws.Range("I4").Value = "Md Biplob Hossain"
$ws.Range("J4").Value = 1839178501
$ws.Range("K4").Value = "Natore"
$ws.Range("L4").Value = "Lalpur"
$ws.Range("M4").Value = "ZSO-0022"
$ws.Range("N4").Value = "Nabinogor, Lalpur, Natore."
$ws.Range("P4").Value = 1839178501
$ws.Range("Q4").Value = "C"
$ws.Range("R4").Value = "Rural"
$ws.Range("S4").Value = "bKash"
$ws.Range("T4").Value = 1839178501

$ws.Range("A5").Value = "DEL-0179"
$ws.Range("B5").Value = "DSR-0349"
$ws.Range("C5").Value = "Nabil Computer"
$ws.Range("D5").Value = "Bagha"
$ws.Range("E5").Value = "Kamruzzaman Pias"
$ws.Range("G5").Value = "GO"
$ws.Range("I5").Value = "Kamruzzaman Pias"
$ws.Range("J5").Value = 1710002659
$ws.Range("K5").Value = "Rajshahi"
$ws.Range("L5").Value = "Bagha"
$ws.Range("M5").Value = "ZSO-0022"
$ws.Range("N5").Value = "Bagha Bazar, Rajshahi."
$ws.Range("P5").Value = 1710002659
$ws.Range("Q5").Value = "C"
$ws.Range("R5").Value = "Rural"
$ws.Range("S5").Value = "bKash"
$ws.Range("T5").Value = 1710002659

$ws.Range("A6").Value = "DEL-0179"
$ws.Range("B6").Value = "DSR-0248"
$ws.Range("C6").Value = "Monia Telecom"
$ws.Range("D6").Value = "Jonail"
$ws.Range("E6").Value = "Md Munjur Alam"
$ws.Range("G6").Value = "GO"
$ws.Range("I6").Value = "Md Munjur Alam"
$ws.Range("J6").Value = 1770800299
$ws.Range("K6").Value = "Natore"
$ws.Range("L6").Value = "Baraigram"
$ws.Range("M6").Value = "ZSO-0022"
$ws.Range("N6").Value = "Jonail Bazar, Baraigram, Natore."
$ws.Range("P6").Value = 1770800299
$ws.Range("Q6").Value = "C"
$ws.Range("R6").Value = "Rural"
$ws.Range("S6").Value = "bKash"
$ws.Range("T6").Value = 1770800299

$ws.Range("A7").Value = "DEL-0179"
$ws.Range("B7").Value = "DSR-0247"
$ws.Range("C7").Value = "Alamin Telecom"
$ws.Range("D7").Value = "Singra"
$ws.Range("E7").Value = "Md Khorshed Vandary"
$ws.Range("G7").Value = "GO"
$ws.Range("I7").Value = "Md Khorshed Vandary"
$ws.Range("J7").Value = 1762676884
$ws.Range("K7").Value = "Natore"
$ws.Range("L7").Value = "Singra"
$ws.Range("M7").Value = "ZSO-0022"
$ws.Range("N7").Value = "Sigra."
$ws.Range("P7").Value = 1762676884
$ws.Range("Q7").Value = "C"
$ws.Range("R7").Value = "Rural"
$ws.Range("S7").Value = "bKash"
$ws.Range("T7").Value = 1762676884


# Column C width change (bestFit was recalculated due to new longer text, width 28.43 -> 36)
$ws.Columns.Item(3).ColumnWidth = 35.1666666666667

# Update active cell selection to C20
$ws.Range("C20").Select()
